$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing row 103, shifting the old rows 103-104
# down to 104-105.
$ws.Rows.Item(103).Insert()

# Populate the new row 103 with the new record (Fecha 2021-09-09 / serial 44448).
$ws.Cells.Item(103, 1).Value = 4
$ws.Cells.Item(103, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(103, 3).Value = "Los Lagos"
$ws.Cells.Item(103, 4).Value = 44448
$ws.Cells.Item(103, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(103, 5).Value = 10
$ws.Cells.Item(103, 6).Value = "Fruta"
$ws.Cells.Item(103, 7).Value = 100101
$ws.Cells.Item(103, 8).Value = "Berries"
$ws.Cells.Item(103, 9).Value = 100101007
$ws.Cells.Item(103, 10).Value = "Kiwi"
$ws.Cells.Item(103, 11).Value = "Hayward"
$ws.Cells.Item(103, 12).Value = "Especial"
$ws.Cells.Item(103, 13).Value = 200
$ws.Cells.Item(103, 14).Value = 20000
$ws.Cells.Item(103, 15).Value = 20000
$ws.Cells.Item(103, 16).Value = 20000
$ws.Cells.Item(103, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(103, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(103, 19).Value = 1333
$ws.Cells.Item(103, 20).Value = 15
